# Running all modules except watchlist
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Watchlist module row (row 6) Runmode from Y to N
$ws.Range("C6").Value = "N"

# Update the active selection to M6 as recorded after the edit
$ws.Activate()
$ws.Range("M6").Select()
